$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact / No display for ContactDetail" row (row 11),
# which collapses the two identical Contact rows into one.
$ws.Rows(11).Delete()

# The remaining former-"Contact" row (row 10) becomes the new "Jurisdiction" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
